# Auto-generated: update cached market-price columns (H:N) across all 8 job sheets
# per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2368  # H19 (was 1278.5714)
$ws.Cells.Item(19, 9).Value = 2324  # I19 (was 1350)
$ws.Cells.Item(19, 10).Value = 2500  # J19 (was 1183.3334)
$ws.Cells.Item(19, 11).Value = 2324  # K19 (was 1350)
$ws.Cells.Item(19, 12).Value = 2500  # L19 (was 1183.3334)
$ws.Cells.Item(19, 13).Value = -2149  # M19 (was -1175)
$ws.Cells.Item(19, 14).Value = -2850  # N19 (was -1533.3334)
$ws.Cells.Item(28, 8).Value = 10989.3  # H28 (was 9453.416999999999)
$ws.Cells.Item(28, 9).Value = 14643.286  # I28 (was 14650.286)
$ws.Cells.Item(28, 10).Value = 2463.3333  # J28 (was 2177.8)
$ws.Cells.Item(28, 11).Value = 14643.286  # K28 (was 14650.286)
$ws.Cells.Item(28, 12).Value = 2463.3333  # L28 (was 2177.8)
$ws.Cells.Item(28, 13).Value = -14158.286  # M28 (was -14165.286)
$ws.Cells.Item(28, 14).Value = -3433.3333  # N28 (was -3147.8)
$ws.Cells.Item(41, 8).Value = 381.1  # H41 (was 401.77777)
$ws.Cells.Item(41, 9).Value = 257.625  # I41 (was 266.57144)
$ws.Cells.Item(41, 11).Value = 257.625  # K41 (was 266.57144)
$ws.Cells.Item(41, 13).Value = 182.375  # M41 (was 173.42856)
$ws.Cells.Item(92, 8).Value = 1446.3077  # H92 (was 1457.2858)
$ws.Cells.Item(92, 9).Value = 1212.5555  # I92 (was 1228.4546)
$ws.Cells.Item(92, 10).Value = 1972.25  # J92 (was 2296.3333)
$ws.Cells.Item(92, 11).Value = 1212.5555  # K92 (was 1228.4546)
$ws.Cells.Item(92, 12).Value = 1972.25  # L92 (was 2296.3333)
$ws.Cells.Item(92, 13).Value = 35.44450000000006  # M92 (was 19.54539999999997)
$ws.Cells.Item(92, 14).Value = -4468.25  # N92 (was -4792.3333)
$ws.Cells.Item(98, 8).Value = 612.375  # H98 (was 582.5294)
$ws.Cells.Item(98, 9).Value = 373.54544  # I98 (was 351.16666)
$ws.Cells.Item(98, 11).Value = 373.54544  # K98 (was 351.16666)
$ws.Cells.Item(98, 13).Value = 1124.45456  # M98 (was 1146.83334)
$ws.Cells.Item(116, 8).Value = 6798  # H116 (was 6664.3335)
$ws.Cells.Item(116, 10).Value = 7333.3335  # J116 (was 6999)
$ws.Cells.Item(116, 12).Value = 7333.3335  # L116 (was 6999)
$ws.Cells.Item(116, 14).Value = -14217.3335  # N116 (was -13883)
$ws.Cells.Item(122, 8).Value = 612.375  # H122 (was 582.5294)
$ws.Cells.Item(122, 9).Value = 373.54544  # I122 (was 351.16666)
$ws.Cells.Item(122, 11).Value = 1120.63632  # K122 (was 1053.49998)
$ws.Cells.Item(122, 13).Value = 1329.36368  # M122 (was 1396.50002)
$ws.Cells.Item(135, 8).Value = 657  # H135 (was 775)
$ws.Cells.Item(135, 10).Value = 1192.5  # J135 (was 1250)
$ws.Cells.Item(135, 12).Value = 10732.5  # L135 (was 11250)
$ws.Cells.Item(135, 14).Value = -15802.5  # N135 (was -16320)
$ws.Cells.Item(141, 8).Value = 3702.8  # H141 (was 4111.3076)
$ws.Cells.Item(141, 9).Value = 3702.8  # I141 (was 4111.3076)
$ws.Cells.Item(141, 11).Value = 11108.4  # K141 (was 12333.9228)
$ws.Cells.Item(141, 13).Value = -5928.400000000001  # M141 (was -7153.9228)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1998.5454  # H2 (was 2018.9)
$ws.Cells.Item(2, 10).Value = 2832.5  # J2 (was 3040)
$ws.Cells.Item(2, 12).Value = 2832.5  # L2 (was 3040)
$ws.Cells.Item(2, 14).Value = -3058.5  # N2 (was -3266)
$ws.Cells.Item(56, 8).Value = 6000  # H56 (was 0)
$ws.Cells.Item(56, 9).Value = 6000  # I56 (was 0)
$ws.Cells.Item(56, 11).Value = 6000  # K56 (was 0)
$ws.Cells.Item(56, 13).Value = -5258  # M56 (was None)
$ws.Cells.Item(61, 8).Value = 3599.2  # H61 (was 3750)
$ws.Cells.Item(61, 9).Value = 3499  # I61 (was 3500)
$ws.Cells.Item(61, 11).Value = 3499  # K61 (was 3500)
$ws.Cells.Item(61, 13).Value = -3287  # M61 (was -3288)
$ws.Cells.Item(74, 8).Value = 5383.4165  # H74 (was 5468.273)
$ws.Cells.Item(74, 9).Value = 4708.7  # I74 (was 4737.4443)
$ws.Cells.Item(74, 11).Value = 4708.7  # K74 (was 4737.4443)
$ws.Cells.Item(74, 13).Value = -3834.7  # M74 (was -3863.4443)
$ws.Cells.Item(77, 8).Value = 5383.4165  # H77 (was 5468.273)
$ws.Cells.Item(77, 9).Value = 4708.7  # I77 (was 4737.4443)
$ws.Cells.Item(77, 11).Value = 23543.5  # K77 (was 23687.2215)
$ws.Cells.Item(77, 13).Value = -19175.5  # M77 (was -19319.2215)
$ws.Cells.Item(92, 8).Value = 44966.668  # H92 (was 36633.332)
$ws.Cells.Item(92, 10).Value = 44966.668  # J92 (was 36633.332)
$ws.Cells.Item(92, 12).Value = 44966.668  # L92 (was 36633.332)
$ws.Cells.Item(92, 14).Value = -49958.668  # N92 (was -41625.332)
$ws.Cells.Item(116, 8).Value = 1998.5454  # H116 (was 2018.9)
$ws.Cells.Item(116, 10).Value = 2832.5  # J116 (was 3040)
$ws.Cells.Item(116, 12).Value = 2832.5  # L116 (was 3040)
$ws.Cells.Item(116, 14).Value = -7420.5  # N116 (was -7628)
$ws.Cells.Item(132, 8).Value = 1625.4333  # H132 (was 1706.2222)
$ws.Cells.Item(132, 9).Value = 1655.3928  # I132 (was 1746.24)
$ws.Cells.Item(132, 11).Value = 4966.178400000001  # K132 (was 5238.72)
$ws.Cells.Item(132, 13).Value = -2436.178400000001  # M132 (was -2708.72)
$ws.Cells.Item(136, 8).Value = 3599.2  # H136 (was 3750)
$ws.Cells.Item(136, 9).Value = 3499  # I136 (was 3500)
$ws.Cells.Item(136, 11).Value = 10497  # K136 (was 10500)
$ws.Cells.Item(136, 13).Value = -7947  # M136 (was -7950)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1998.5454  # H3 (was 2018.9)
$ws.Cells.Item(3, 10).Value = 2832.5  # J3 (was 3040)
$ws.Cells.Item(3, 12).Value = 2832.5  # L3 (was 3040)
$ws.Cells.Item(3, 14).Value = -3060.5  # N3 (was -3268)
$ws.Cells.Item(82, 8).Value = 33525.7  # H82 (was 31377.908)
$ws.Cells.Item(82, 9).Value = 18419  # I82 (was 16289.25)
$ws.Cells.Item(82, 11).Value = 18419  # K82 (was 16289.25)
$ws.Cells.Item(82, 13).Value = -18036  # M82 (was -15906.25)
$ws.Cells.Item(85, 8).Value = 33525.7  # H85 (was 31377.908)
$ws.Cells.Item(85, 9).Value = 18419  # I85 (was 16289.25)
$ws.Cells.Item(85, 11).Value = 18419  # K85 (was 16289.25)
$ws.Cells.Item(85, 13).Value = -17093  # M85 (was -14963.25)
$ws.Cells.Item(105, 8).Value = 2163.6667  # H105 (was 1816.25)
$ws.Cells.Item(105, 9).Value = 1573  # I105 (was 1217.8889)
$ws.Cells.Item(105, 11).Value = 1573  # K105 (was 1217.8889)
$ws.Cells.Item(105, 13).Value = 174  # M105 (was 529.1111000000001)
$ws.Cells.Item(107, 8).Value = 884.1111  # H107 (was 934.75)
$ws.Cells.Item(107, 9).Value = 744.625  # I107 (was 782.5714)
$ws.Cells.Item(107, 11).Value = 744.625  # K107 (was 782.5714)
$ws.Cells.Item(107, 13).Value = 1175.375  # M107 (was 1137.4286)
$ws.Cells.Item(134, 8).Value = 3457.16  # H134 (was 3517.96)
$ws.Cells.Item(134, 9).Value = 3501.25  # I134 (was 3564.5833)
$ws.Cells.Item(134, 11).Value = 10503.75  # K134 (was 10693.7499)
$ws.Cells.Item(134, 13).Value = -7968.75  # M134 (was -8158.749899999999)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2588.5898  # H31 (was 2497.6829)
$ws.Cells.Item(31, 9).Value = 2104.1155  # I31 (was 2005.6072)
$ws.Cells.Item(31, 11).Value = 2104.1155  # K31 (was 2005.6072)
$ws.Cells.Item(31, 13).Value = -1809.1155  # M31 (was -1710.6072)
$ws.Cells.Item(34, 8).Value = 2588.5898  # H34 (was 2497.6829)
$ws.Cells.Item(34, 9).Value = 2104.1155  # I34 (was 2005.6072)
$ws.Cells.Item(34, 11).Value = 2104.1155  # K34 (was 2005.6072)
$ws.Cells.Item(34, 13).Value = -1902.1155  # M34 (was -1803.6072)
$ws.Cells.Item(86, 8).Value = 14166.5  # H86 (was 15000)
$ws.Cells.Item(86, 9).Value = 13333  # I86 (was 0)
$ws.Cells.Item(86, 11).Value = 13333  # K86 (was 0)
$ws.Cells.Item(86, 13).Value = -12210  # M86 (was None)
$ws.Cells.Item(89, 8).Value = 14166.5  # H89 (was 15000)
$ws.Cells.Item(89, 9).Value = 13333  # I89 (was 0)
$ws.Cells.Item(89, 11).Value = 66665  # K89 (was 0)
$ws.Cells.Item(89, 13).Value = -61049  # M89 (was None)
$ws.Cells.Item(99, 8).Value = 5002.7334  # H99 (was 5504.273)
$ws.Cells.Item(99, 9).Value = 5054.3  # I99 (was 5506)
$ws.Cells.Item(99, 10).Value = 4899.6  # J99 (was 5499.6665)
$ws.Cells.Item(99, 11).Value = 5054.3  # K99 (was 5506)
$ws.Cells.Item(99, 12).Value = 4899.6  # L99 (was 5499.6665)
$ws.Cells.Item(99, 13).Value = -3556.3  # M99 (was -4008)
$ws.Cells.Item(99, 14).Value = -7895.6  # N99 (was -8495.666499999999)
$ws.Cells.Item(126, 8).Value = 5002.7334  # H126 (was 5504.273)
$ws.Cells.Item(126, 9).Value = 5054.3  # I126 (was 5506)
$ws.Cells.Item(126, 10).Value = 4899.6  # J126 (was 5499.6665)
$ws.Cells.Item(126, 11).Value = 15162.9  # K126 (was 16518)
$ws.Cells.Item(126, 12).Value = 14698.8  # L126 (was 16498.9995)
$ws.Cells.Item(126, 13).Value = -12692.9  # M126 (was -14048)
$ws.Cells.Item(126, 14).Value = -19638.8  # N126 (was -21438.9995)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1718.8148  # H4 (was 1454.125)
$ws.Cells.Item(4, 9).Value = 1818.6875  # I4 (was 1391.5714)
$ws.Cells.Item(4, 11).Value = 5456.0625  # K4 (was 4174.7142)
$ws.Cells.Item(4, 13).Value = -5344.0625  # M4 (was -4062.7142)
$ws.Cells.Item(5, 8).Value = 3040  # H5 (was 2760)
$ws.Cells.Item(5, 9).Value = 4264.3335  # I5 (was 3398.25)
$ws.Cells.Item(5, 11).Value = 12793.0005  # K5 (was 10194.75)
$ws.Cells.Item(5, 13).Value = -12681.0005  # M5 (was -10082.75)
$ws.Cells.Item(94, 8).Value = 977.7778  # H94 (was 988.8889)
$ws.Cells.Item(94, 9).Value = 977.7778  # I94 (was 988.8889)
$ws.Cells.Item(94, 11).Value = 2933.3334  # K94 (was 2966.6667)
$ws.Cells.Item(94, 13).Value = -2257.3334  # M94 (was -2290.6667)
$ws.Cells.Item(135, 8).Value = 3040  # H135 (was 2760)
$ws.Cells.Item(135, 9).Value = 4264.3335  # I135 (was 3398.25)
$ws.Cells.Item(135, 11).Value = 38379.0015  # K135 (was 30584.25)
$ws.Cells.Item(135, 13).Value = -35844.0015  # M135 (was -28049.25)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2999.6667  # H102 (was 2184.2)
$ws.Cells.Item(102, 9).Value = 2999.5  # I102 (was 1980.25)
$ws.Cells.Item(102, 11).Value = 2999.5  # K102 (was 1980.25)
$ws.Cells.Item(102, 13).Value = -1377.5  # M102 (was -358.25)
$ws.Cells.Item(132, 8).Value = 9432.666999999999  # H132 (was 10149)
$ws.Cells.Item(132, 9).Value = 8999.5  # I132 (was 9999)
$ws.Cells.Item(132, 11).Value = 26998.5  # K132 (was 29997)
$ws.Cells.Item(132, 13).Value = -24468.5  # M132 (was -27467)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2000  # H7 (was 0)
$ws.Cells.Item(7, 10).Value = 2000  # J7 (was 0)
$ws.Cells.Item(7, 12).Value = 2000  # L7 (was 0)
$ws.Cells.Item(7, 14).Value = -2224  # N7 (was None)
$ws.Cells.Item(40, 8).Value = 4999.5  # H40 (was 0)
$ws.Cells.Item(40, 9).Value = 4999.5  # I40 (was 0)
$ws.Cells.Item(40, 11).Value = 4999.5  # K40 (was 0)
$ws.Cells.Item(40, 13).Value = -4863.5  # M40 (was None)
$ws.Cells.Item(46, 8).Value = 4230.6924  # H46 (was 4454.5454)
$ws.Cells.Item(46, 9).Value = 2000  # I46 (was 0)
$ws.Cells.Item(46, 10).Value = 4416.5835  # J46 (was 4454.5454)
$ws.Cells.Item(46, 11).Value = 2000  # K46 (was 0)
$ws.Cells.Item(46, 12).Value = 4416.5835  # L46 (was 4454.5454)
$ws.Cells.Item(46, 13).Value = -1812  # M46 (was None)
$ws.Cells.Item(46, 14).Value = -4792.5835  # N46 (was -4830.5454)
$ws.Cells.Item(55, 8).Value = 307.73685  # H55 (was 305.6842)
$ws.Cells.Item(55, 9).Value = 265.6  # I55 (was 259.63635)
$ws.Cells.Item(55, 10).Value = 354.55554  # J55 (was 369)
$ws.Cells.Item(55, 11).Value = 265.6  # K55 (was 259.63635)
$ws.Cells.Item(55, 12).Value = 354.55554  # L55 (was 369)
$ws.Cells.Item(55, 13).Value = -92.60000000000002  # M55 (was -86.63634999999999)
$ws.Cells.Item(55, 14).Value = -700.5555400000001  # N55 (was -715)
$ws.Cells.Item(100, 8).Value = 1350.75  # H100 (was 1320.6)
$ws.Cells.Item(100, 10).Value = 0  # J100 (was 1200)
$ws.Cells.Item(100, 12).Value = 0  # L100 (was 1200)
$ws.Cells.Item(100, 14).Value = $null  # N100 (was -2282)
$ws.Cells.Item(126, 8).Value = 2000  # H126 (was 0)
$ws.Cells.Item(126, 10).Value = 2000  # J126 (was 0)
$ws.Cells.Item(126, 12).Value = 6000  # L126 (was 0)
$ws.Cells.Item(126, 14).Value = -10940  # N126 (was None)
$ws.Cells.Item(132, 8).Value = 5701.647  # H132 (was 6026.125)
$ws.Cells.Item(132, 9).Value = 5442.8  # I132 (was 5740.8)
$ws.Cells.Item(132, 10).Value = 5809.5  # J132 (was 6155.8184)
$ws.Cells.Item(132, 11).Value = 16328.4  # K132 (was 17222.4)
$ws.Cells.Item(132, 12).Value = 17428.5  # L132 (was 18467.4552)
$ws.Cells.Item(132, 13).Value = -13798.4  # M132 (was -14692.4)
$ws.Cells.Item(132, 14).Value = -22488.5  # N132 (was -23527.4552)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(108, 8).Value = 0  # H108 (was 55000)
$ws.Cells.Item(108, 10).Value = 0  # J108 (was 55000)
$ws.Cells.Item(108, 12).Value = 0  # L108 (was 55000)
$ws.Cells.Item(108, 14).Value = $null  # N108 (was -62680)
$ws.Cells.Item(113, 8).Value = 7453  # H113 (was 7949.7144)
$ws.Cells.Item(113, 9).Value = 17375  # I113 (was 20750.2)
$ws.Cells.Item(113, 11).Value = 52125  # K113 (was 62250.60000000001)
$ws.Cells.Item(113, 13).Value = -49955  # M113 (was -60080.60000000001)
$ws.Cells.Item(122, 8).Value = 1242.8572  # H122 (was 1599.5)
$ws.Cells.Item(122, 9).Value = 1232.8334  # I122 (was 1599.5)
$ws.Cells.Item(122, 10).Value = 1303  # J122 (was 0)
$ws.Cells.Item(122, 11).Value = 3698.5002  # K122 (was 4798.5)
$ws.Cells.Item(122, 12).Value = 3909  # L122 (was 0)
$ws.Cells.Item(122, 13).Value = -1248.5002  # M122 (was -2348.5)
$ws.Cells.Item(122, 14).Value = -8809  # N122 (was None)
$ws.Cells.Item(132, 8).Value = 1691.95  # H132 (was 1767.8422)
$ws.Cells.Item(132, 9).Value = 1535.7222  # I132 (was 1611.3529)
$ws.Cells.Item(132, 11).Value = 4607.1666  # K132 (was 4834.0587)
$ws.Cells.Item(132, 13).Value = -2077.1666  # M132 (was -2304.0587)
